$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "hierarchical naïve (MARIO)"
$ws.Range("B5").Value = "satura la memoria, tempo infinito"

$ws.Range("B5").Select()
